$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: boAt Airdopes 131/138 product, price Rs. 1,299
$ws.Range("A1").Value = "boAt Airdopes 131/138 Twin Wireless Earbuds with IWP Technology, Bluetooth V5.0, Immersive Audio, Up to 15H Total Playback, Instant Voice Assistant and Type-C Charging,Bluetooth Earphone (Active Black)"
$ws.Range("B1").Value = "Rs. 1,299"

# Row 2: NBOX MARATHON, reworded description, price Rs. 749
$ws.Range("A2").Value = "NBOX MARATHON Over Ear Bluetooth Neckband 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Silver"
$ws.Range("B2").Value = "Rs. 749"

# Row 3: new NBOX Air1 TWS product, price Rs. 749
$ws.Range("A3").Value = "NBOX Air1 TWS On Ear True Wireless (TWS) 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Black"
$ws.Range("B3").Value = "Rs. 749"

# Row 4: NBOX STAR 40 HOURS product, price Rs. 700
$ws.Range("A4").Value = "NBOX STAR 40 HOURS MUSIC PLAYBACK WIRELESS NECKBAND WITH DOLBY EFFECT BASS SOUND IPX5 WITH MASSIVE MUSIC PLAYBACK,BLUETOOTH HEADPHONE,BLUETOOTH EARPHONE FOR NBOX"
$ws.Range("B4").Value = "Rs. 700"

# Row 5: boAt Airdopes 121v2 product, price Rs. 1,299
$ws.Range("A5").Value = "boAt Airdopes 121v2 On Ear True Wireless (TWS) 14 Hours Playback IPX7(Water Resistant) Active Noise cancellation -Bluetooth V 5.0 Black"
$ws.Range("B5").Value = "Rs. 1,299"
